# Insert a new data row at row 190 (pushing the existing rows 190-299 down
# to 191-300, carrying their formatting/styles with them - this also grows
# the sheet's used range from A1:R299 to A1:R300), then populate the newly
# inserted row with the new weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(190).Insert()

$ws.Cells.Item(190, 1).Value  = 9
$ws.Cells.Item(190, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(190, 3).Value  = "Metropolitana"
$ws.Cells.Item(190, 4).Value  = 44777
$ws.Cells.Item(190, 5).Value  = 13
$ws.Cells.Item(190, 6).Value  = 300000001
$ws.Cells.Item(190, 7).Value  = "Rabanito"
$ws.Cells.Item(190, 8).Value  = "Sin especificar"
$ws.Cells.Item(190, 9).Value  = "Primera"
$ws.Cells.Item(190, 10).Value = 7000
$ws.Cells.Item(190, 11).Value = 2500
$ws.Cells.Item(190, 12).Value = 3000
$ws.Cells.Item(190, 13).Value = 2750
$ws.Cells.Item(190, 14).Value = "$/cien unidades (volumen en unidades)"
$ws.Cells.Item(190, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(190, 16).Value = 28
$ws.Cells.Item(190, 17).Value = 100
$ws.Cells.Item(190, 18).Value = "Hortaliza"
